# Update "想去人数" (F column) counts for the 展览 and 全部类型 sheets.
# Both sheets carry the same underlying listing data, so the same
# row -> new-value updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 137
    5  = 93
    6  = 131
    7  = 1273
    8  = 1540
    10 = 398
    12 = 154
    14 = 64
    15 = 108
    17 = 308
    18 = 325
    19 = 1737
    23 = 670
    25 = 337
    26 = 4194
    28 = 271
    29 = 1091
    30 = 487
    32 = 554
    33 = 23
    34 = 256
    36 = 140
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
